# Update the "Correspond Handoff Datetime" (E) and "Correspond Handback
# DateTime" (H) columns for the 9b2c1858... row (row 3) on both the
# "zh-cn" and "de-de" status sheets, reflecting a newer handback run.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-21 22:46:03"
$wsZhCn.Range("H3").Value = "2016-03-21 22:46:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-21 22:46:10"
$wsDeDe.Range("H3").Value = "2016-03-21 22:46:40"
